$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "51.629.48"
$ws.Cells.Item(2, 5).Value = "  +1.04%  "

$ws.Cells.Item(3, 4).Value = "3.040.01"
$ws.Cells.Item(3, 5).Value = "  +2.65%  "

$ws.Cells.Item(4, 5).Value = "  +0.04%  "

$ws.Cells.Item(5, 4).Value = "384.98"
$ws.Cells.Item(5, 5).Value = "  +1.28%  "

$ws.Cells.Item(6, 4).Value = "102.91"
$ws.Cells.Item(6, 5).Value = "  +0.53%  "

$ws.Cells.Item(7, 5).Value = "  -0.08%  "

$ws.Cells.Item(8, 5).Value = "  +0.00%  "

$ws.Cells.Item(9, 4).Value = "0.588"
$ws.Cells.Item(9, 5).Value = "  -0.38%  "

$ws.Cells.Item(10, 4).Value = "36.94"
$ws.Cells.Item(10, 5).Value = "  +0.98%  "

$ws.Cells.Item(11, 5).Value = "  +0.08%  "

$ws.Cells.Item(12, 4).Value = "0.0862"
$ws.Cells.Item(12, 5).Value = "  +0.93%  "

$ws.Cells.Item(13, 4).Value = "3.509.65"
$ws.Cells.Item(13, 5).Value = "  +2.46%  "

$ws.Cells.Item(14, 4).Value = "18.70"
$ws.Cells.Item(14, 5).Value = "  +2.00%  "

$ws.Cells.Item(15, 4).Value = "7.78"
$ws.Cells.Item(15, 5).Value = "  +0.01%  "

$ws.Cells.Item(16, 4).Value = "3.043.15"
$ws.Cells.Item(16, 5).Value = "  +2.44%  "

$ws.Cells.Item(17, 4).Value = "0.976"
$ws.Cells.Item(17, 5).Value = "  -2.37%  "

$ws.Cells.Item(18, 4).Value = "10.58"
$ws.Cells.Item(18, 5).Value = "  -11.86%  "

$ws.Cells.Item(19, 4).Value = "51.648.77"
$ws.Cells.Item(19, 5).Value = "  +0.91%  "

$ws.Cells.Item(20, 4).Value = "3.09"
$ws.Cells.Item(20, 5).Value = "  -0.36%  "

$ws.Cells.Item(21, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(21, 4).Value = "12.38"
$ws.Cells.Item(21, 5).Value = "  -0.22%  "

$ws.Cells.Item(22, 2).Value = "ShibaInu"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(22, 4).Value = "0.0₃0963"
$ws.Cells.Item(22, 5).Value = "  +0.18%  "

$ws.Cells.Item(23, 4).Value = "69.91"
$ws.Cells.Item(23, 5).Value = "  -0.34%  "

$ws.Cells.Item(24, 4).Value = "267.16"
$ws.Cells.Item(24, 5).Value = "  -0.29%  "

$ws.Cells.Item(25, 4).Value = "3.18"
$ws.Cells.Item(25, 5).Value = "  -2.84%  "

$ws.Cells.Item(26, 5).Value = "  +5.29%  "

$ws.Cells.Item(27, 4).Value = "7.43"
$ws.Cells.Item(27, 5).Value = "  +2.77%  "

$ws.Cells.Item(28, 5).Value = "  +4.11%  "

$ws.Cells.Item(29, 4).Value = "26.37"
$ws.Cells.Item(29, 5).Value = "  +1.81%  "

$ws.Cells.Item(30, 5).Value = "  -0.03%  "

$ws.Cells.Item(31, 5).Value = "  -1.83%  "

$ws.Cells.Item(32, 4).Value = "10.28"
$ws.Cells.Item(32, 5).Value = "  -1.55%  "

$ws.Cells.Item(33, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(33, 4).Value = "34.13"
$ws.Cells.Item(33, 5).Value = "  -1.05%  "

$ws.Cells.Item(34, 2).Value = "Toncoin"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(34, 4).Value = "2.07"
$ws.Cells.Item(34, 5).Value = "  -0.14%  "

$ws.Cells.Item(35, 4).Value = "50.54"
$ws.Cells.Item(35, 5).Value = "  -0.93%  "

$ws.Cells.Item(36, 4).Value = "0.0449"
$ws.Cells.Item(36, 5).Value = "  +2.72%  "

$ws.Cells.Item(37, 5).Value = "  -0.16%  "

$ws.Cells.Item(38, 4).Value = "3.38"
$ws.Cells.Item(38, 5).Value = "  +3.88%  "

$ws.Cells.Item(39, 2).Value = "TheGraph"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(39, 4).Value = "0.286"
$ws.Cells.Item(39, 5).Value = "  +6.36%  "

$ws.Cells.Item(40, 2).Value = "Celestia"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(40, 4).Value = "17.01"
$ws.Cells.Item(40, 5).Value = "  +2.45%  "

$ws.Cells.Item(41, 5).Value = "  +1.49%  "

$ws.Cells.Item(42, 5).Value = "  -0.30%  "

$ws.Cells.Item(43, 4).Value = "127.36"
$ws.Cells.Item(43, 5).Value = "  +2.05%  "

$ws.Cells.Item(44, 4).Value = "2.53"
$ws.Cells.Item(44, 5).Value = "  +0.78%  "

$ws.Cells.Item(45, 5).Value = "  +4.06%  "

$ws.Cells.Item(46, 4).Value = "21.64"
$ws.Cells.Item(46, 5).Value = "  -0.10%  "

$ws.Cells.Item(47, 4).Value = "2.48"
$ws.Cells.Item(47, 5).Value = "  +4.33%  "

$ws.Cells.Item(48, 5).Value = "  +4.12%  "

$ws.Cells.Item(49, 4).Value = "2.034.51"
$ws.Cells.Item(49, 5).Value = "  -0.81%  "

$ws.Cells.Item(50, 4).Value = "3.336.35"
$ws.Cells.Item(50, 5).Value = "  +2.56%  "

$ws.Cells.Item(51, 4).Value = "0.209"
$ws.Cells.Item(51, 5).Value = "  +7.70%  "

$wb.Save()